$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(5, 9).Value = "sd"
$ws.Cells.Item(5, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(12, 9).Value = "sd"
$ws.Cells.Item(12, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(18, 9).Value = "sv"
$ws.Cells.Item(18, 10).Value = "Statement-opinion"
$ws.Cells.Item(32, 9).Value = "sv"
$ws.Cells.Item(32, 10).Value = "Statement-opinion"
$ws.Cells.Item(33, 9).Value = "sd"
$ws.Cells.Item(33, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(42, 9).Value = "sd"
$ws.Cells.Item(42, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(44, 9).Value = "sd"
$ws.Cells.Item(44, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(51, 9).Value = "aa"
$ws.Cells.Item(51, 10).Value = "Agree/Accept"
$ws.Cells.Item(55, 9).Value = "b"
$ws.Cells.Item(55, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(59, 9).Value = "sd"
$ws.Cells.Item(59, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(78, 9).Value = "sd"
$ws.Cells.Item(78, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(93, 9).Value = "b"
$ws.Cells.Item(93, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(100, 9).Value = "sd"
$ws.Cells.Item(100, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(105, 9).Value = "sv"
$ws.Cells.Item(105, 10).Value = "Statement-opinion"
$ws.Cells.Item(109, 9).Value = "b"
$ws.Cells.Item(109, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(124, 9).Value = "sd"
$ws.Cells.Item(124, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(127, 9).Value = "ba"
$ws.Cells.Item(127, 10).Value = "Appreciation"
$ws.Cells.Item(134, 9).Value = "aa"
$ws.Cells.Item(134, 10).Value = "Agree/Accept"
$ws.Cells.Item(139, 9).Value = "aa"
$ws.Cells.Item(139, 10).Value = "Agree/Accept"
$ws.Cells.Item(143, 9).Value = "sv"
$ws.Cells.Item(143, 10).Value = "Statement-opinion"
$ws.Cells.Item(150, 9).Value = "b"
$ws.Cells.Item(150, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(165, 9).Value = "%"
$ws.Cells.Item(165, 10).Value = "Uninterpretable"
$ws.Cells.Item(168, 9).Value = "%"
$ws.Cells.Item(168, 10).Value = "Uninterpretable"
$ws.Cells.Item(175, 9).Value = "sd"
$ws.Cells.Item(175, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(183, 9).Value = "sd"
$ws.Cells.Item(183, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(189, 9).Value = "sd"
$ws.Cells.Item(189, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(192, 9).Value = "sv"
$ws.Cells.Item(192, 10).Value = "Statement-opinion"
$ws.Cells.Item(198, 9).Value = "sd"
$ws.Cells.Item(198, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(215, 9).Value = "b"
$ws.Cells.Item(215, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(218, 9).Value = "sd"
$ws.Cells.Item(218, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(224, 9).Value = "aa"
$ws.Cells.Item(224, 10).Value = "Agree/Accept"
$ws.Cells.Item(225, 9).Value = "sd"
$ws.Cells.Item(225, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(226, 9).Value = "sd"
$ws.Cells.Item(226, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(228, 9).Value = "sd"
$ws.Cells.Item(228, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(230, 9).Value = "aa"
$ws.Cells.Item(230, 10).Value = "Agree/Accept"
$ws.Cells.Item(231, 9).Value = "aa"
$ws.Cells.Item(231, 10).Value = "Agree/Accept"
$ws.Cells.Item(232, 9).Value = "b"
$ws.Cells.Item(232, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(233, 9).Value = "%"
$ws.Cells.Item(233, 10).Value = "Uninterpretable"
$ws.Cells.Item(234, 9).Value = "%"
$ws.Cells.Item(234, 10).Value = "Uninterpretable"
$ws.Cells.Item(235, 9).Value = "%"
$ws.Cells.Item(235, 10).Value = "Uninterpretable"
$ws.Cells.Item(257, 9).Value = "b"
$ws.Cells.Item(257, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(265, 9).Value = "sv"
$ws.Cells.Item(265, 10).Value = "Statement-opinion"
$ws.Cells.Item(278, 9).Value = "sd"
$ws.Cells.Item(278, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(279, 9).Value = "aa"
$ws.Cells.Item(279, 10).Value = "Agree/Accept"
$ws.Cells.Item(285, 9).Value = "sd"
$ws.Cells.Item(285, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(305, 9).Value = "aa"
$ws.Cells.Item(305, 10).Value = "Agree/Accept"
$ws.Cells.Item(309, 9).Value = "aa"
$ws.Cells.Item(309, 10).Value = "Agree/Accept"
$ws.Cells.Item(316, 9).Value = "aa"
$ws.Cells.Item(316, 10).Value = "Agree/Accept"
$ws.Cells.Item(318, 9).Value = "sd"
$ws.Cells.Item(318, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(320, 9).Value = "aa"
$ws.Cells.Item(320, 10).Value = "Agree/Accept"
$ws.Cells.Item(328, 9).Value = "b"
$ws.Cells.Item(328, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(342, 9).Value = "b"
$ws.Cells.Item(342, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(366, 9).Value = "aa"
$ws.Cells.Item(366, 10).Value = "Agree/Accept"
$ws.Cells.Item(384, 9).Value = "aa"
$ws.Cells.Item(384, 10).Value = "Agree/Accept"
$ws.Cells.Item(388, 9).Value = "sd"
$ws.Cells.Item(388, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(404, 9).Value = "sd"
$ws.Cells.Item(404, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(417, 9).Value = "sd"
$ws.Cells.Item(417, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(426, 9).Value = "sd"
$ws.Cells.Item(426, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(432, 9).Value = "aa"
$ws.Cells.Item(432, 10).Value = "Agree/Accept"
$ws.Cells.Item(435, 9).Value = "sd"
$ws.Cells.Item(435, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(436, 9).Value = "aa"
$ws.Cells.Item(436, 10).Value = "Agree/Accept"
$ws.Cells.Item(439, 9).Value = "sd"
$ws.Cells.Item(439, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(447, 9).Value = "sd"
$ws.Cells.Item(447, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(448, 9).Value = "sd"
$ws.Cells.Item(448, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(451, 9).Value = "sd"
$ws.Cells.Item(451, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(456, 9).Value = "sd"
$ws.Cells.Item(456, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(457, 9).Value = "sd"
$ws.Cells.Item(457, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(458, 9).Value = "b"
$ws.Cells.Item(458, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(467, 9).Value = "sv"
$ws.Cells.Item(467, 10).Value = "Statement-opinion"
$ws.Cells.Item(471, 9).Value = "aa"
$ws.Cells.Item(471, 10).Value = "Agree/Accept"
$ws.Cells.Item(474, 9).Value = "sd"
$ws.Cells.Item(474, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(475, 9).Value = "aa"
$ws.Cells.Item(475, 10).Value = "Agree/Accept"
$ws.Cells.Item(477, 9).Value = "sd"
$ws.Cells.Item(477, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(481, 9).Value = "b"
$ws.Cells.Item(481, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(483, 9).Value = "sd"
$ws.Cells.Item(483, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(493, 9).Value = "b"
$ws.Cells.Item(493, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(496, 9).Value = "sd"
$ws.Cells.Item(496, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(508, 9).Value = "aa"
$ws.Cells.Item(508, 10).Value = "Agree/Accept"
$ws.Cells.Item(511, 9).Value = "sv"
$ws.Cells.Item(511, 10).Value = "Statement-opinion"
$ws.Cells.Item(514, 9).Value = "%"
$ws.Cells.Item(514, 10).Value = "Uninterpretable"
$ws.Cells.Item(524, 9).Value = "sd"
$ws.Cells.Item(524, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(526, 9).Value = "sd"
$ws.Cells.Item(526, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(527, 9).Value = "aa"
$ws.Cells.Item(527, 10).Value = "Agree/Accept"
$ws.Cells.Item(528, 9).Value = "%"
$ws.Cells.Item(528, 10).Value = "Uninterpretable"
$ws.Cells.Item(531, 9).Value = "sv"
$ws.Cells.Item(531, 10).Value = "Statement-opinion"
$ws.Cells.Item(534, 9).Value = "sv"
$ws.Cells.Item(534, 10).Value = "Statement-opinion"
$ws.Cells.Item(545, 9).Value = "aa"
$ws.Cells.Item(545, 10).Value = "Agree/Accept"
$ws.Cells.Item(548, 9).Value = "b"
$ws.Cells.Item(548, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(554, 9).Value = "aa"
$ws.Cells.Item(554, 10).Value = "Agree/Accept"
$ws.Cells.Item(557, 9).Value = "sd"
$ws.Cells.Item(557, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(558, 9).Value = "aa"
$ws.Cells.Item(558, 10).Value = "Agree/Accept"
$ws.Cells.Item(574, 9).Value = "sd"
$ws.Cells.Item(574, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(576, 9).Value = "sd"
$ws.Cells.Item(576, 10).Value = "Statement-non-opinion"
